$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.700.07"
$ws.Range("E2").Value = "  +0.01%  "

$ws.Range("D3").Value = "2.468.34"
$ws.Range("E3").Value = "  +0.00%  "

$ws.Range("E4").Value = "  +0.05%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "320.62"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +1.38%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "92.20"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.87%  "

$ws.Range("E7").Value = "  -0.22%  "


$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.507"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -1.39%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "32.90"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +0.20%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0854"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +1.19%  "

$ws.Range("E12").Value = "  -0.84%  "

$ws.Range("D13").Value = "2.849.61"
$ws.Range("E13").Value = "  +0.02%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "6.88"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -0.30%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "15.49"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -1.56%  "

$ws.Range("D16").Value = "2.463.66"
$ws.Range("E16").Value = "  -0.10%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.789"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +0.82%  "

$ws.Range("D18").Value = "41.616.40"
$ws.Range("E18").Value = "  -0.12%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "6.44"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -0.96%  "

$ws.Range("E20").Value = "  -1.35%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "71.85"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +0.96%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "11.20"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -2.50%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "239.80"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +0.15%  "

$ws.Range("E24").Value = "  +0.90%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "1.94"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +0.27%  "

$ws.Range("E26").Value = "  +0.01%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "24.80"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +0.26%  "

$ws.Range("E28").Value = "  -1.88%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "9.70"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -0.97%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "36.08"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +1.40%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "155.35"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -0.44%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "5.42"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -1.51%  "

$ws.Range("E33").Value = "  -0.04%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.0763"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -0.19%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "2.57"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -0.52%  "

$ws.Range("E36").Value = "  -2.86%  "

$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "2.92"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +1.15%  "

$ws.Range("B38").Value = "Stellar"
$ws.Range("C38").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.116"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +1.10%  "

$ws.Range("E39").Value = "  +1.86%  "

$ws.Range("E40").Value = "  -0.42%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "3.98"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -0.15%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "2.32"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -7.54%  "

$ws.Range("D43").Value = "2.002.53"
$ws.Range("E43").Value = "  +1.46%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.0282"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -1.02%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "18.56"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -2.14%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "2.95"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +0.39%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "9.46"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +4.50%  "

$ws.Range("D48").Value = "2.730.92"
$ws.Range("E48").Value = "  +1.08%  "

$ws.Range("B49").Value = "BitcoinSV"
$ws.Range("C49").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "75.76"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +3.83%  "

$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "97.14"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -0.11%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "66.90"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -0.14%  "

